$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New IMEI entries added for today's sales (16.11.19)
$ws.Range("A7").Value = 357653104891481
$ws.Range("B7").Value = "B65"
$ws.Rows.Item(7).RowHeight = 20.1

$ws.Range("A8").Value = 359998100004623
$ws.Range("B8").Value = "D40"
$ws.Rows.Item(8).RowHeight = 20.1

$ws.Range("A13").Value = 355580101577908
$ws.Range("B13").Value = "R40"
$ws.Rows.Item(13).RowHeight = 20.1

$ws.Range("A17").Value = 358444100265204
$ws.Range("B17").Value = "s40"
$ws.Rows.Item(17).RowHeight = 20.1

# Update the view to where the user was last working
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F14").Select()
